$d = $word.ActiveDocument

# 1) Merge the three runs "  introduce your research " / "questions and hypothesis." / "]"
#    into a single run by re-issuing the combined text via Find & Replace (collapses runs).
$d.Content.Find.Execute(
    " introduce your research questions and hypothesis.]", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " introduce your research questions and hypothesis.]", 2) | Out-Null

# 2) Merge the five runs that make up the "introduce the research methods..." sentence.
$d.Content.Find.Execute(
    " introduce the research methods and data sources you used for the analysis]", $true, $false, $false,
    $false, $false, $true, 1, $false,
    " introduce the research methods and data sources you used for the analysis]", 2) | Out-Null

# 3) Merge the closing smart-quote + parenthesis runs after "Speaking".
#    (Restrict Find to a tight Range so the differently-formatted "Speaking"
#    run in front is left untouched and only the two ") "-ish runs combine.)
$rightQuote = [char]8221
$probe = $d.Content
$probe.Find.Execute("Speaking", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $probe.End
$sub = $d.Range($pos, $pos + 2)
$sub.Find.Execute($rightQuote + ")", $true, $false, $false, $false, $false, $true, 1, $false,
    $rightQuote + ")", 2) | Out-Null

# 4) Merge the closing smart-quote + parenthesis/space runs after "Taller Children".
$probe = $d.Content
$probe.Find.Execute("Taller Children", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$pos = $probe.End
$sub = $d.Range($pos, $pos + 3)
$sub.Find.Execute($rightQuote + ") ", $true, $false, $false, $false, $false, $true, 1, $false,
    $rightQuote + ") ", 2) | Out-Null

# 5) Drop the stray leading hyphen before "Results" (this is the only
#    hyphenated "-Results" run in the document).
$r = $d.Content
$r.Find.Execute("-Results ", $true, $false, $false, $false, $false, $true, 1, $false,
    "Results ", 2) | Out-Null

# 6) Move the "_GoBack" bookmark from the end of the "Narrative, p93" paragraph to
#    right after the ")" that closes "(power of narrative//storytelling)" in the
#    "-Discussion" bullet immediately below the "Results" bullet just edited
#    (there is an earlier, unrelated occurrence of the same phrase near the top
#    of the document, so resume the search from where we already are).
$d.Bookmarks("_GoBack").Delete()

$r.Collapse(0)
$r.Find.Execute("power of narrative//storytelling)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r.Collapse(0)
$d.Bookmarks.Add("_GoBack", $r) | Out-Null
